# Auto update: 2025-12-06 21:20:02
# Swap the AIG / MetLife rows (name + ticker) and refresh the numeric metrics
# for all four holdings (rows 2-5) to the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 becomes "American International Group, I" / "AIG" ---
$ws.Range("B3").Value = "American International Group, I"
$ws.Range("C3").Value = "AIG"

# --- Row 4 becomes "MetLife, Inc." / "MET" ---
$ws.Range("B4").Value = "MetLife, Inc."
$ws.Range("C4").Value = "MET"

# --- Row 2 (UNH) refreshed metrics ---
$ws.Range("D2").Value = 330.91
$ws.Range("E2").Value = 56.7
$ws.Range("F2").Value = 0.35
$ws.Range("G2").Value = 60
$ws.Range("K2").Value = 58.9
$ws.Range("N2").Value = 52.28493729186943

# --- Row 3 (AIG) refreshed metrics ---
$ws.Range("D3").Value = 77.03
$ws.Range("E3").Value = 44.9
$ws.Range("F3").Value = 1.14
$ws.Range("H3").Value = 46
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 53.7
$ws.Range("N3").Value = 52.28493729186943

# --- Row 4 (MetLife) refreshed metrics ---
$ws.Range("D4").Value = 78.66
$ws.Range("E4").Value = 50.1
$ws.Range("F4").Value = 2.74
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 23
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 53.7
$ws.Range("N4").Value = 52.28493729186943

# --- Row 5 (Prudential) refreshed metrics ---
$ws.Range("D5").Value = 111.68
$ws.Range("E5").Value = 70.2
$ws.Range("F5").Value = 3.17
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 33
$ws.Range("K5").Value = 40.9
$ws.Range("N5").Value = 52.28493729186943
